$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- Update the "ulais11xx" identifiers in column A with typo'd values ---
# (Column B already holds the "urbsNN" values and is left untouched.)
# NOTE: rows 5/6/7 are intentionally written out of order (7,6,5) so that the
# shared-string table ends up with the same index assignment as the source
# workbook (string table order follows first-use order).
$ws2.Range("A1").Value = "yulais1145"
$ws2.Range("A2").Value = "yulais1146"
$ws2.Range("A3").Value = "uylais1147"
$ws2.Range("A4").Value = "uylais1148"
$ws2.Range("A7").Value = "ulais1y151"
$ws2.Range("A6").Value = "ulaisy1150"
$ws2.Range("A5").Value = "ulayis1149"
$ws2.Range("A8").Value = "ulyais1152"

# --- Update the selected cell on the testCitizen sheet ---
$ws2.Activate()
$ws2.Range("A8").Select() | Out-Null

# --- Widen column A slightly ---
$ws2.Columns.Item(1).ColumnWidth = 19.5
